# Data_output. Corrigiendo función plot
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ypfd")

# --- Fix existing rows 91 and 92 ---
$ws.Cells.Item(91, 6).Value = 57441218.95
$ws.Cells.Item(91, 7).Value = 326
$ws.Cells.Item(91, 8).Value = 561

$ws.Cells.Item(92, 6).Value = 109329325.7
$ws.Cells.Item(92, 8).Value = 1326

# --- Append new rows 575-585 ---
$newRows = @(
    @("2024-05-14", 25074.8, 24599, 25180, 24206, 5374027820.1, 216554, 4531),
    @("2024-05-15", 26287, 25165, 26400, 25165, 6915901084.35, 264300, 5209),
    @("2024-05-16", 26060, 26215, 26352.15, 25793, 6186295653.05, 237033, 4212),
    @("2024-05-17", 25900, 26000, 26269.95, 25810.25, 6320768245.3, 242778, 4042),
    @("2024-05-20", 27200, 25900, 27200, 25672.85, 6509051725.95, 245452, 4406),
    @("2024-05-21", 28218.45, 27141.5, 28352.35, 27100, 9171119839.799999, 331953, 5755),
    @("2024-05-22", 27408.2, 28150, 28218.45, 27285.95, 7986168929.4, 288649, 7024),
    @("2024-05-23", 27000, 27411, 28150, 26845.55, 5339281026, 194893, 5964),
    @("2024-05-24", 26120, 27100, 27385.85, 25897.75, 8617342873.4, 324692, 6297),
    @("2024-05-27", 26840, 26500, 26990, 26487.05, 580551257.95, 4, 1384),
    @("2024-05-28", 27669, 27420.4, 27700, 26655, 9398066727.5, 344887, 5556)
)

$startRow = 575
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    # Column A holds a date-like string; force text storage (no special style)
    # so it matches the plain inline/shared string used for every other date cell.
    $dateCell = $ws.Cells.Item($r, 1)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $row[0]
    $dateCell.Style = "Normal"

    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
}
